# Update to countries & provincias Spain (9 Sep 2020 refresh, ~22:10)
#
# The source COVID-19 stats feed refreshed totals for several countries;
# because the "Pais" sheet is kept sorted descending by "Casos totales"
# (column B), a handful of countries swapped adjacent rows as a result:
#   - Cabo Verde overtook Ruanda and Surinam (rows 120-122)
#   - Republica de Chipre overtook Togo (rows 155-156)
#   - Santa Lucia overtook Nueva Caledonia (rows 204-205)
#   - Montserrat overtook Islas Malvinas (rows 214-215)
# Estados Unidos, Sudafrica, Alemania, Israel, Angola, Islas Feroe and
# Monaco simply got refreshed figures without changing rank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados a ..." timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 9 de Septiembre de 2020 a las 22:10"

# Rows that changed: row number, Pais, Casos totales, Nuevos casos, Casos
# activas, Recuperados, Casos criticos, Muertes hoy, Muertes
$data = @(
  @(4, "Estados Unidos", 6537006, 22775, 3821394, 2520832, 0, 750, 194780),
  @(11, "Sudafrica", 642431, 1990, 569935, 57328, 0, 82, 15168),
  @(24, "Alemania", 256334, 1378, 230600, 16323, 0, 2, 9411),
  @(28, "Israel", 141097, 3532, 108354, 31689, 0, 14, 1054),
  @(120, "Cabo Verde", 4473, 73, 3915, 515, 0, 1, 43),
  @(121, "Ruanda", 4439, 0, 2307, 2112, 0, 0, 20),
  @(122, "Surinam", 4419, 0, 3595, 733, 0, 0, 91),
  @(134, "Angola", 3092, 59, 1245, 1721, 0, 2, 126),
  @(155, "Republica de Chipre", 1514, 3, 1237, 255, 0, 0, 22),
  @(156, "Togo", 1513, 0, 1127, 352, 0, 0, 34),
  @(179, "Islas Feroe", 415, 1, 410, 5, 0, 0, 0),
  @(191, "Monaco", 161, 5, 107, 53, 0, 0, 1),
  @(204, "Santa Lucia", 27, 1, 26, 1, 0, 0, 0),
  @(205, "Nueva Caledonia", 26, 0, 25, 1, 0, 0, 0),
  @(214, "Montserrat", 13, 0, 12, 0, 0, 0, 1),
  @(215, "Islas Malvinas", 13, 0, 13, 0, 0, 0, 0)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
    $ws.Cells.Item($r, 3).Value = $item[3]
    $ws.Cells.Item($r, 4).Value = $item[4]
    $ws.Cells.Item($r, 5).Value = $item[5]
    $ws.Cells.Item($r, 6).Value = $item[6]
    $ws.Cells.Item($r, 7).Value = $item[7]
    $ws.Cells.Item($r, 8).Value = $item[8]
}
